$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Affiliation corrections (shared-string content fixes) ---
# Row 106: Guillaume Joubert -- affiliation changed from Heudiasyc to IMT Atlantique
$ws.Range("B106").Value = "IMT Atlantique"
# Row 165: Niels Agatz -- affiliation changed from Delft University of Technology to Rotterdam School of Management
$ws.Range("B165").Value = "Rotterdam School of Management"

# --- New committee members appended as rows 183-210 ---
$ws.Range("A183").Value = "Lei Zhao"
$ws.Range("B183").Value = "Tsinghua University"
$ws.Range("C183").Value = "China"
$ws.Range("A184").Value = "Angelica Del Rocio Lozano Cuevas"
$ws.Range("B184").Value = "Universidad Nacional Autónoma de México`t"
$ws.Range("C184").Value = "Mexico"
$ws.Range("A185").Value = "Panagiotis Repoussis"
$ws.Range("B185").Value = "Athens University of Economics and Business"
$ws.Range("C185").Value = "Greece"
$ws.Range("A186").Value = "Alexandre Jacquillat"
$ws.Range("B186").Value = "Massachusetts Institute of Technology"
$ws.Range("C186").Value = "USA"
$ws.Range("A187").Value = "Anne Goodchild"
$ws.Range("B187").Value = "University of Washington"
$ws.Range("C187").Value = "USA"
$ws.Range("A188").Value = "Ann Melissa Campbell"
$ws.Range("B188").Value = "University of Iowa"
$ws.Range("C188").Value = "USA"
$ws.Range("A189").Value = "Stefan Minner"
$ws.Range("B189").Value = "Technical University of Munich`t"
$ws.Range("C189").Value = "Germany"
$ws.Range("A190").Value = "Mihaela Popa"
$ws.Range("B190").Value = "University Politehnica of Bucharest"
$ws.Range("C190").Value = "Romania"
$ws.Range("A191").Value = "Ruibin Bai"
$ws.Range("B191").Value = "University of Nottingham Ningbo"
$ws.Range("C191").Value = "China"
$ws.Range("A192").Value = "Richard Wong"
$ws.Range("A193").Value = "Anant Balakrishnan"
$ws.Range("B193").Value = "University of Texas"
$ws.Range("C193").Value = "USA"
$ws.Range("A194").Value = "Antonio Mauttone"
$ws.Range("B194").Value = "Universidad de la República"
$ws.Range("C194").Value = "Uruguay"
$ws.Range("A195").Value = "Claudio Sterle"
$ws.Range("B195").Value = "University of Naples"
$ws.Range("C195").Value = "Italy"
$ws.Range("A196").Value = "Margaretha Gansterer"
$ws.Range("B196").Value = "University of Vienna`t"
$ws.Range("C196").Value = "Austria"
$ws.Range("A197").Value = "Karl F. Doerner"
$ws.Range("B197").Value = "University of Vienna"
$ws.Range("C197").Value = "Austria"
$ws.Range("A198").Value = "Leandro Callegari Coelho"
$ws.Range("B198").Value = "Laval University"
$ws.Range("C198").Value = "Canada"
$ws.Range("A199").Value = "Nadia Lahrichi"
$ws.Range("B199").Value = "Polytechnique Montréal"
$ws.Range("C199").Value = "Canada"
$ws.Range("A200").Value = "Kjetil Fagerholt"
$ws.Range("B200").Value = "Norwegian University of Science and Technology"
$ws.Range("C200").Value = "Norway"
$ws.Range("A201").Value = "Jorge Mendoza Gimenez"
$ws.Range("B201").Value = "HEC Montréal"
$ws.Range("C201").Value = "Canada"
$ws.Range("A202").Value = "Fausto Errico"
$ws.Range("B202").Value = "CIRRELT and École de technologie supérieure de Montréal"
$ws.Range("C202").Value = "Canada"
$ws.Range("A203").Value = "Guido Perboli"
$ws.Range("B203").Value = "Politecnico di Torino, Italy and CIRRELT, Canada"
$ws.Range("A204").Value = "Federico Malucelli"
$ws.Range("B204").Value = "Politecnico di Milano"
$ws.Range("C204").Value = "Italy"
$ws.Range("A205").Value = "Massimo Di Francesco"
$ws.Range("B205").Value = "Università di Cagliari"
$ws.Range("C205").Value = "Italy"
$ws.Range("A206").Value = "Manuel Iori"
$ws.Range("B206").Value = "University of Modena and Reggio Emilia"
$ws.Range("C206").Value = "Italy"
$ws.Range("A207").Value = "Antonio Frangioni"
$ws.Range("B207").Value = "Università di Pisa"
$ws.Range("C207").Value = "Italy"
$ws.Range("A208").Value = "Maria Elena Bruni"
$ws.Range("B208").Value = "University of Calabria"
$ws.Range("C208").Value = "Italy"
$ws.Range("A209").Value = "Gianfranco Guastaroba"
$ws.Range("B209").Value = "University of Brescia"
$ws.Range("C209").Value = "Italy"
$ws.Range("A210").Value = "Bilge Atasoy"
$ws.Range("B210").Value = "Delft University of Technology"
$ws.Range("C210").Value = "Netherlands"

# --- Restore active selection to match the saved view state ---
$ws.Range("D189").Select()
